# Apply updated "dSF" (column F) values to specific rows on Sheet1.
# This mirrors a data "repull" where the dSF column is refreshed with
# newly recalculated figures while the original dS0 (column E) snapshot
# is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = -4
    5  = 5
    7  = 1
    8  = 2
    10 = 1
    12 = -1
    14 = 0
    15 = -2
    19 = -1
    21 = 2
    29 = -2
    30 = 2
    34 = -3
    35 = -7
    37 = -2
    38 = -10
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
